# Improve form customer:
#  - Penambahan pilihan jenis customer (semua, umum dan reseller)
#  - Penambahan fitur pencarian nama customer
#
# Concretely (per the target OOXML diff): insert a new "DISKON RESELLER"
# header column into the "customer" sheet, between "PLAFON PIUTANG"'s old
# position (I) and the rest of the header row, shifting "PLAFON PIUTANG"
# from column I to column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer")

# Insert a new column at I (shifts the existing "PLAFON PIUTANG" column,
# and its column width/style definition, from I to J).
$ws.Columns.Item(9).Insert()

# New header cell I1: same header (blue fill / white bold text) formatting
# as the rest of the header row, with the new label.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 9).Value = "DISKON RESELLER"

# Give the new column its own best-fit style/width (matches the other
# data columns' style index) instead of the default column formatting.
$ws.Columns.Item(9).ColumnWidth = 16

# Update the active selection as left by the editor.
$ws.Range("E4").Select()
